# Update Work Week and Social Spending
# Refresh the Kenya GDP per Capita series (Data sheet, column E) with new
# values, and extend the series with six more years (2011-2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated GDP-per-capita series: 61 values for the existing rows (years
# 1950-2010, rows 2-62) followed by 6 values for the newly appended rows
# (years 2011-2016, rows 63-68).
$vals = @("1038","1229","1063","1009","1095","1144","1173","1176","1156","1148","1157","1093","1117","1138","1208","1184","1294","1317","1366","1404","1458","1500","1524","1546","1564","1502","1503","1580","1656","1651","1675","1647","1680","1627","1594","1604","1658","1698","1741","1761","1780","1777.78378054864","1724.8744312691","1698.25414622437","1729.71716115145","1789.89454405253","1851.28138215243","1846.3738357664","1889.97310366846","1920.83394493269","1914.8863918102","1972.38110614001","1961.89493984161","1997.77682002573","2066.21019811963","2157.84156954474","2255.25856515838","2380.77633261763","2357.23254112996","2404.62957991628","2579.53091765354","2712","2765","2854","2942","3046","3169")

# Rows 2-62 already exist: just overwrite column E with the refreshed data.
for ($i = 0; $i -lt 61; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = "'" + $vals[$i]
}

# Rows 63-68 are brand new: years 2011-2016, same Country Code/Name/Indicator
# as the rest of the series, with the last 6 refreshed values.
$years = @(2011, 2012, 2013, 2014, 2015, 2016)
for ($i = 0; $i -lt 6; $i++) {
    $row = 63 + $i
    $ws.Cells.Item($row, 1).Value = 404
    $ws.Cells.Item($row, 2).Value = "Kenya"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $years[$i]
    $ws.Cells.Item($row, 5).Value = "'" + $vals[61 + $i]
}
